$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE changes from "001" to "004" (keep as text, leading zeros matter)
$ws.Range("J2").Value = "'004"

# N2: REPORT_DATE text value
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Numeric financial figures for row 2
$ws.Range("O2").Value = 73964805.98999999
$ws.Range("P2").Value = 419294553.5
$ws.Range("Q2").Value = 341673363.44
$ws.Range("R2").Value = 26.9368250212
$ws.Range("S2").Value = 250886563.97
$ws.Range("T2").Value = 250886563.97
$ws.Range("U2").Value = 30.3328387271
$ws.Range("V2").Value = 31224572.59
$ws.Range("W2").Value = 19123494.68
$ws.Range("X2").Value = 1751668.4
$ws.Range("Y2").Value = 80624470.20999999
$ws.Range("Z2").Value = 82960097.61
$ws.Range("AA2").Value = 8995291.619999999
$ws.Range("AG2").Value = 2191714.64
$ws.Range("AP2").Value = 29.8732504331
$ws.Range("AQ2").Value = 42.45921690275
$ws.Range("AR2").Value = 42.42
$ws.Range("AS2").Value = 70736440.70999999
$ws.Range("AT2").Value = 40.127216881439
